$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates: ID 4 -> 3, Temperatura 35.7 -> 37.2,
# Fecha de ingreso 2020-12-05 -> 2020-12-04, Hora de ingreso 02:38:13.203773 -> 02:45:47.400847
$ws.Range("A2").Value = 3
$ws.Range("F2").Value = 37.2
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "2020-12-04"
$ws.Range("H2").Value = "02:45:47.400847"

# Row 3 updates: ID 5 -> 4, Nro documento 1193474912 -> 1339998889, Uso del tapabocas No -> Si,
# Temperatura 36 -> 35.7, Hora de ingreso 02:40:18.030912 -> 02:38:13.203773,
# Ingreso Denegado -> Aceptado, Nombres Isabela -> Anuel, Apellidos Acevedo Garcia -> AA
$ws.Range("A3").Value = 4
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1339998889"
$ws.Range("E3").Value = "Si"
$ws.Range("F3").Value = 35.7
$ws.Range("H3").Value = "02:38:13.203773"
$ws.Range("I3").Value = "Aceptado"
$ws.Range("J3").Value = "Anuel"
$ws.Range("K3").Value = "AA"
